$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.653.11"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.564.23"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.45"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.13"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.55"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.81"
$ws.Range("E13").Value = "  -2.84%  "
$ws.Range("D14").Value = "3.018.66"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "62.539.27"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").Value = "2.555.48"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.11"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.86"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.55"
$ws.Range("E24").Value = "  +4.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.58"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.91"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.16"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("E30").Value = "  -3.58%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0792"
$ws.Range("E31").Value = "  -4.13%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "451.37"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.35"
$ws.Range("E33").Value = "  -0.21%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.63"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.394"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.69"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "39.67"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "155.72"
$ws.Range("E42").Value = "  +3.05%  "
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.630"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.79"
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.39"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -4.67%  "
